# Apply the edit described by the diff:
# - Sheet "2" becomes the active/selected tab (activeTab=1 in workbook, tabSelected on sheet2, removed from sheet1)
# - Selection on sheet 1 stays N13 (unchanged)
# - Selection on sheet 2 changes from M18 to L7
# - Sheet "2": G4 gets value "H" (matching shared string used elsewhere, style-less like B4..H4)
# - Sheet "2": I4's value "H" is cleared (kept style s=5)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("1")
$ws2 = $wb.Worksheets.Item("2")

# Fill in the missing G4 cell on sheet "2" with "H"
$ws2.Range("G4").Value = "H"

# Clear the value in I4 on sheet "2" (keep cell/style, remove content)
$ws2.Range("I4").ClearContents()

# Update selections
$ws1.Range("N13").Select()
$ws2.Range("L7").Select()

# Make sheet "2" the active sheet/tab
$ws2.Activate()
